$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Row: Deep - Mid
Replace-Text "0.005" "-0.006"
Replace-Text "0.430" "-0.348"
Replace-Text "0.904" "0.935"

# Row: Deep - Shallow
Replace-Text "0.003" "-0.012"
Replace-Text "0.241" "-0.680"
Replace-Text "0.968" "0.777"

# Row: Mid - Shallow
Replace-Text "-0.002" "-0.006"
Replace-Text "-0.188" "-0.332"
Replace-Text "0.981" "0.941"

# Columns shared across all three data rows: std.error (0.011 -> 0.018), df (24 -> 22)
Replace-Text "0.011" "0.018"
Replace-Text "24" "22"
